$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.606816530227661
$ws.Range("B1").Value = 1.272009968757629
$ws.Range("C1").Value = 1.680297017097473
$ws.Range("D1").Value = 2.54114556312561
$ws.Range("E1").Value = 6.473733901977539
